$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Julio de 2020 a las 19:24"

# Row 4
$ws.Range("B4").Value = 4051681
$ws.Range("C4").Value = 23112
$ws.Range("D4").Value = 1896310
$ws.Range("E4").Value = 2010042
$ws.Range("G4").Value = 376
$ws.Range("H4").Value = 145329

# Row 5
$ws.Range("B5").Value = 2178159
$ws.Range("C5").Value = 11627
$ws.Range("E5").Value = 630361
$ws.Range("G5").Value = 231
$ws.Range("H5").Value = 81828

# Row 6
$ws.Range("B6").Value = 1238013
$ws.Range("C6").Value = 43928
$ws.Range("D6").Value = 782780
$ws.Range("E6").Value = 425346
$ws.Range("G6").Value = 1117
$ws.Range("H6").Value = 29887

# Row 11
$ws.Range("B11").Value = 336402
$ws.Range("C11").Value = 1719
$ws.Range("E11").Value = 18439

# Row 12
$ws.Range("B12").Value = 314631
$ws.Range("C12").Value = 1357
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 28426

# Row 18
$ws.Range("B18").Value = 222402
$ws.Range("C18").Value = 902
$ws.Range("D18").Value = 205214
$ws.Range("E18").Value = 11643
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 5545

# Row 21
$ws.Range("B21").Value = 204298
$ws.Range("C21").Value = 408
$ws.Range("E21").Value = 6517
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 9181

# Row 24
$ws.Range("B24").Value = 112005
$ws.Range("C24").Value = 308
$ws.Range("D24").Value = 98042
$ws.Range("E24").Value = 5095
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 8868

# Row 42
$ws.Range("A42").Value = "Israel"
$ws.Range("B42").Value = 55695
$ws.Range("C42").Value = 1653
$ws.Range("D42").Value = 23205
$ws.Range("E42").Value = 32060
$ws.Range("G42").Value = 5
$ws.Range("H42").Value = 430

# Row 43
$ws.Range("A43").Value = "Panama"
$ws.Range("B43").Value = 55153
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 30075
$ws.Range("E43").Value = 23919
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 1159

# Row 59
$ws.Range("B59").Value = 25819
$ws.Range("C59").Value = 17
$ws.Range("E59").Value = 701
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 1754

# Row 61
$ws.Range("B61").Value = 24872
$ws.Range("C61").Value = 594
$ws.Range("D61").Value = 16983
$ws.Range("E61").Value = 6778
$ws.Range("G61").Value = 11
$ws.Range("H61").Value = 1111

# Row 67
$ws.Range("B67").Value = 17962
$ws.Range("C67").Value = 220
$ws.Range("D67").Value = 15636
$ws.Range("E67").Value = 2041
$ws.Range("G67").Value = 5
$ws.Range("H67").Value = 285

# Row 80
$ws.Range("B80").Value = 9547
$ws.Range("C80").Value = 135
$ws.Range("D80").Value = 5071
$ws.Range("E80").Value = 4034
$ws.Range("G80").Value = 10
$ws.Range("H80").Value = 442

# Row 102
$ws.Range("B102").Value = 4077
$ws.Range("C102").Value = 29
$ws.Range("E102").Value = 2503
$ws.Range("G102").Value = 3
$ws.Range("H102").Value = 200

# Row 109
$ws.Range("A109").Value = "Libano"
$ws.Range("B109").Value = 3104
$ws.Range("C109").Value = 124
$ws.Range("D109").Value = 1582
$ws.Range("E109").Value = 1479
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 43

# Row 110
$ws.Range("B110").Value = 3103
$ws.Range("C110").Value = 59
$ws.Range("E110").Value = 691

# Row 111
$ws.Range("A111").Value = "Guinea Ecuatorial"
$ws.Range("B111").Value = 3071
$ws.Range("D111").Value = 842
$ws.Range("E111").Value = 2178
$ws.Range("H111").Value = 51

# Row 113
$ws.Range("B113").Value = 2839
$ws.Range("C113").Value = 15
$ws.Range("E113").Value = 187

# Row 114
$ws.Range("B114").Value = 2745
$ws.Range("C114").Value = 15
$ws.Range("E114").Value = 670

# Row 119
$ws.Range("A119").Value = "Libia"
$ws.Range("B119").Value = 2176
$ws.Range("C119").Value = 88
$ws.Range("D119").Value = 489
$ws.Range("E119").Value = 1634
$ws.Range("G119").Value = 3
$ws.Range("H119").Value = 53

# Row 120
$ws.Range("B120").Value = 2154
$ws.Range("C120").Value = 47
$ws.Range("D120").Value = 1132
$ws.Range("E120").Value = 1001

# Row 121
$ws.Range("A121").Value = "Hong Kong"
$ws.Range("B121").Value = 2132
$ws.Range("C121").Value = 112
$ws.Range("D121").Value = 1344
$ws.Range("E121").Value = 774
$ws.Range("H121").Value = 14

# Row 130
$ws.Range("B130").Value = 1731
$ws.Range("C130").Value = 4
$ws.Range("D130").Value = 1288
$ws.Range("E130").Value = 377

# Row 148
$ws.Range("A148").Value = "Principado de Andorra"
$ws.Range("C148").Value = 5
$ws.Range("E148").Value = 34
$ws.Range("H148").Value = 52

# Row 149
$ws.Range("A149").Value = "Republica del Chad"
$ws.Range("B149").Value = 889
$ws.Range("D149").Value = 805
$ws.Range("E149").Value = 9
$ws.Range("H149").Value = 75

# Row 157
$ws.Range("B157").Value = 645
$ws.Range("C157").Value = 6
$ws.Range("D157").Value = 558
$ws.Range("E157").Value = 84

# Row 189
$ws.Range("B189").Value = 82
$ws.Range("C189").Value = 1
$ws.Range("E189").Value = 58

# Row 210
$ws.Range("A210").Value = "Islas Malvinas"

# Row 211
$ws.Range("A211").Value = "Groenlandia"
